$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 645
$ws.Range("F5").Value = 2883
$ws.Range("F9").Value = 288
$ws.Range("F10").Value = 6708
$ws.Range("F11").Value = 24
$ws.Range("F12").Value = 304
$ws.Range("F14").Value = 973
$ws.Range("F15").Value = 581
$ws.Range("F16").Value = 1442
$ws.Range("F18").Value = 1088
$ws.Range("F19").Value = 2167
$ws.Range("F20").Value = 1410
$ws.Range("F21").Value = 632
$ws.Range("F22").Value = 83
$ws.Range("F23").Value = 1055
$ws.Range("F24").Value = 71
$ws.Range("F25").Value = 150
$ws.Range("F26").Value = 297
$ws.Range("F27").Value = 1604
$ws.Range("F28").Value = 769
$ws.Range("F30").Value = 19
$ws.Range("F31").Value = 21
$ws.Range("F32").Value = 1640
$ws.Range("F33").Value = 1147
$ws.Range("F34").Value = 127
$ws.Range("F37").Value = 366
$ws.Range("F39").Value = 2654
$ws.Range("F40").Value = 66
$ws.Range("F48").Value = 123

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F7").Value = 125
$ws.Range("F14").Value = 46
$ws.Range("F22").Value = 445
$ws.Range("F25").Value = 20
$ws.Range("F35").Value = 27

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 531
$ws.Range("F6").Value = 1719
$ws.Range("F7").Value = 1615
$ws.Range("F9").Value = 2659
$ws.Range("F10").Value = 964
$ws.Range("F11").Value = 844
$ws.Range("F13").Value = 197

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 531
$ws.Range("F4").Value = 1719
$ws.Range("F5").Value = 645
$ws.Range("F6").Value = 2883
$ws.Range("F8").Value = 1615
$ws.Range("F9").Value = 288
$ws.Range("F10").Value = 2659
$ws.Range("F11").Value = 6708
$ws.Range("F12").Value = 964
$ws.Range("F13").Value = 844
$ws.Range("F14").Value = 24
$ws.Range("F15").Value = 304
$ws.Range("F16").Value = 197
$ws.Range("F17").Value = 581
$ws.Range("F18").Value = 1442
$ws.Range("F19").Value = 1088
$ws.Range("F20").Value = 2167
$ws.Range("F21").Value = 1410
$ws.Range("F22").Value = 632
$ws.Range("F23").Value = 83
$ws.Range("F25").Value = 1055
$ws.Range("F26").Value = 71
$ws.Range("F27").Value = 297
$ws.Range("F29").Value = 1604
$ws.Range("F31").Value = 21
$ws.Range("F32").Value = 1640
$ws.Range("F33").Value = 1147
$ws.Range("F37").Value = 445
$ws.Range("F38").Value = 366
$ws.Range("F42").Value = 2654
$ws.Range("F49").Value = 27
